$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap contents of row 81 and row 82 (columns F:V only) ---
$row81vals = $ws.Range("F81:V81").Value2
$row82vals = $ws.Range("F82:V82").Value2
$ws.Range("F81:V81").Value2 = $row82vals
$ws.Range("F82:V82").Value2 = $row81vals

# --- Append new row 92 with match data ---
# Copy formatting from row 91 for the styled cells (A and E columns)
$ws.Range("A91").Copy()
$ws.Range("A92").PasteSpecial(-4122)
$ws.Range("E91").Copy()
$ws.Range("E92").PasteSpecial(-4122)

$ws.Range("A92").Value2 = 91
$ws.Range("B92").Value2 = "denmark"
$ws.Range("C92").Value2 = "superliga"
$ws.Range("D92").Value2 = "2023-2024"
$ws.Range("E92").Value2 = 45254.79166666666
$ws.Range("F92").Value2 = "Hvidovre IF"
$ws.Range("G92").Value2 = 1
$ws.Range("H92").Value2 = "Vejle"
$ws.Range("I92").Value2 = 1
$ws.Range("J92").Value2 = 3.12
$ws.Range("K92").Value2 = "13/11/2023 10:41"
$ws.Range("L92").Value2 = 3.11
$ws.Range("M92").Value2 = "24/11/2023 18:59"
$ws.Range("N92").Value2 = 3.52
$ws.Range("O92").Value2 = "13/11/2023 10:41"
$ws.Range("P92").Value2 = 3.45
$ws.Range("Q92").Value2 = "24/11/2023 18:59"
$ws.Range("R92").Value2 = 2.33
$ws.Range("S92").Value2 = "13/11/2023 10:41"
$ws.Range("T92").Value2 = 2.37
$ws.Range("U92").Value2 = "24/11/2023 18:59"
$ws.Range("V92").Value2 = "https://www.betexplorer.com/football/denmark/superliga/hvidovre-if-vejle/jsLvO5zG/"
